$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.192.31"
$ws.Range("E2").Value = "  -3.10%  "

$ws.Range("D3").Value = "2.480.66"
$ws.Range("E3").Value = "  -3.36%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'564.30"
$ws.Range("E5").Value = "  -3.22%  "

$ws.Range("D6").Value = "'163.44"
$ws.Range("E6").Value = "  -5.16%  "

$ws.Range("E8").Value = "  -1.77%  "

$ws.Range("D9").Value = "2.479.48"
$ws.Range("E9").Value = "  -3.42%  "

$ws.Range("E10").Value = "  -5.89%  "

$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("E12").Value = "  -2.86%  "

$ws.Range("D13").Value = "'4.88"
$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("D14").Value = "2.937.33"
$ws.Range("E14").Value = "  -3.41%  "

$ws.Range("D15").Value = "69.034.84"
$ws.Range("E15").Value = "  -3.22%  "

$ws.Range("E16").Value = "  -3.39%  "

$ws.Range("D17").Value = "'24.15"
$ws.Range("E17").Value = "  -5.32%  "

$ws.Range("D18").Value = "2.473.74"
$ws.Range("E18").Value = "  -3.91%  "

$ws.Range("D19").Value = "'11.10"
$ws.Range("E19").Value = "  -4.94%  "

$ws.Range("E20").Value = "  -7.76%  "

$ws.Range("D21").Value = "'344.68"
$ws.Range("E21").Value = "  -3.82%  "

$ws.Range("E22").Value = "  -3.40%  "

$ws.Range("D23").Value = "'1.90"
$ws.Range("E23").Value = "  -7.63%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "'69.33"
$ws.Range("E25").Value = "  -1.92%  "

$ws.Range("D26").Value = "'3.85"
$ws.Range("E26").Value = "  -6.32%  "

$ws.Range("D28").Value = "'8.61"
$ws.Range("E28").Value = "  -6.10%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "0.0₃0866"
$ws.Range("E30").Value = "  -6.57%  "

$ws.Range("D31").Value = "'7.64"
$ws.Range("E31").Value = "  -4.38%  "

$ws.Range("D32").Value = "'439.96"
$ws.Range("E32").Value = "  -7.67%  "

$ws.Range("E33").Value = "  -8.61%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Value = "'1.69"
$ws.Range("E35").Value = "  -4.81%  "

$ws.Range("D36").Value = "'154.84"
$ws.Range("E36").Value = "  -1.61%  "

$ws.Range("D37").Value = "'0.112"
$ws.Range("E37").Value = "  -5.95%  "

$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("D39").Value = "'18.05"
$ws.Range("E39").Value = "  -4.46%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("E41").Value = "  -3.57%  "

$ws.Range("D42").Value = "'4.56"
$ws.Range("E42").Value = "  -7.23%  "

$ws.Range("E43").Value = "  -4.22%  "

$ws.Range("E44").Value = "  -2.34%  "

$ws.Range("E45").Value = "  -10.14%  "

$ws.Range("E46").Value = "  -9.83%  "

$ws.Range("D47").Value = "'138.00"
$ws.Range("E47").Value = "  -5.78%  "

$ws.Range("E48").Value = "  -4.31%  "

$ws.Range("D49").Value = "'0.509"
$ws.Range("E49").Value = "  -6.11%  "

$ws.Range("E50").Value = "  -2.57%  "

$ws.Range("D51").Value = "'0.570"
$ws.Range("E51").Value = "  -3.27%  "
